$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.001.38'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '2.375.03'
$ws.Range('E3').Value = '  +2.24%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.98'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.22'
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').Value = '  -0.61%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.498'
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.18'
$ws.Range('E10').Value = '  -1.35%  '
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('E12').Value = '  +2.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.21'
$ws.Range('E13').Value = '  -4.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.76'
$ws.Range('E14').Value = '  -0.43%  '
$ws.Range('D15').Value = '2.745.66'
$ws.Range('E15').Value = '  +2.15%  '
$ws.Range('D16').Value = '2.380.25'
$ws.Range('E16').Value = '  +2.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.803'
$ws.Range('E17').Value = '  +1.36%  '
$ws.Range('D18').Value = '42.961.36'
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.19'
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('E20').Value = '  +2.21%  '
$ws.Range('E21').Value = '  -0.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.11'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.99'
$ws.Range('E23').Value = '  -0.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.23'
$ws.Range('E24').Value = '  -1.60%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.44'
$ws.Range('E25').Value = '  +0.62%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.87'
$ws.Range('E27').Value = '  +1.82%  '
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('E29').Value = '  +1.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.53'
$ws.Range('E30').Value = '  -2.82%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('E32').Value = '  +1.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.52'
$ws.Range('E33').Value = '  -1.76%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0733'
$ws.Range('E34').Value = '  +4.38%  '
$ws.Range('E35').Value = '  +5.25%  '
$ws.Range('E36').Value = '  +3.32%  '
$ws.Range('E37').Value = '  -2.99%  '
$ws.Range('E38').Value = '  -0.76%  '
$ws.Range('E39').Value = '  +1.52%  '
$ws.Range('E40').Value = '  +5.27%  '
$ws.Range('E41').Value = '  -0.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '116.76'
$ws.Range('E42').Value = '  -29.77%  '
$ws.Range('D43').Value = '1.943.93'
$ws.Range('E43').Value = '  +0.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0280'
$ws.Range('E44').Value = '  +0.23%  '
$ws.Range('E45').Value = '  +1.86%  '
$ws.Range('E46').Value = '  -1.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.17'
$ws.Range('E47').Value = '  -10.86%  '
$ws.Range('D48').Value = '2.603.30'
$ws.Range('E48').Value = '  +1.88%  '
$ws.Range('E49').Value = '  +1.76%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.03'
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.02'
$ws.Range('E51').Value = '  -2.85%  '
